$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1, matching the style of the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Timestamp values for rows 2-19 (column F), as plain text (matching other data cells' formatting)
$timeTaken = @(
    "2021-10-05 13:40:21.911242",
    "2021-10-05 13:40:21.911255",
    "2021-10-05 13:40:21.911259",
    "2021-10-05 13:40:21.911262",
    "2021-10-05 13:40:21.911266",
    "2021-10-05 13:40:21.911269",
    "2021-10-05 13:40:21.911272",
    "2021-10-05 13:40:21.911275",
    "2021-10-05 13:40:21.911278",
    "2021-10-05 13:40:21.911281",
    "2021-10-05 13:40:21.911284",
    "2021-10-05 13:40:21.911287",
    "2021-10-05 13:40:21.911290",
    "2021-10-05 13:40:21.911293",
    "2021-10-05 13:40:21.911296",
    "2021-10-05 13:40:21.911298",
    "2021-10-05 13:40:21.911302",
    "2021-10-05 13:40:21.911305"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
